$d = $word.ActiveDocument

# 1. Merge the split "${" + "cloneValue" + "}" runs (separated by spellcheck
#    proofErr markers) back into a single run "${cloneValue}".
$d.Content.Find.Execute('${cloneValue}', $false, $false, $false, $false, $false, `
    $true, 1, $false, '${cloneValue}', 2) | Out-Null

# 2. Merge the paragraph-spacing sentence that was split around a
#    grammar-check proofErr marker ("spacing") back into a single run.
$text = " can also format paragraph such as this justified, 12pt before and 12pt after with 1.5 lines spacing paragraph. This formatting can be applied inline or using predefined style as we use to do in Word."
$d.Content.Find.Execute($text, $false, $false, $false, $false, $false, `
    $true, 1, $false, $text, 2) | Out-Null
